$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.109.50'
$ws.Range("E2").Value = '  +2.62%  '
$ws.Range("D3").Value = '1.677.37'
$ws.Range("E3").Value = '  +3.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.74'
$ws.Range("E5").Value = '  +1.81%  '
$ws.Range("E6").Value = '  +2.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  +3.43%  '
$ws.Range("E9").Value = '  +2.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.27'
$ws.Range("E10").Value = '  +5.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0888'
$ws.Range("E11").Value = '  +4.79%  '
$ws.Range("D12").Value = '1.916.36'
$ws.Range("D13").Value = '1.680.60'
$ws.Range("E13").Value = '  +3.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.11'
$ws.Range("E14").Value = '  +2.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.524'
$ws.Range("E15").Value = '  +3.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.11'
$ws.Range("E16").Value = '  +3.74%  '
$ws.Range("D17").Value = '27.128.92'
$ws.Range("E17").Value = '  +2.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '239.19'
$ws.Range("E18").Value = '  +1.14%  '
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.77'
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("E22").Value = '  +4.62%  '
$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.24'
$ws.Range("E23").Value = '  +2.65%  '
$ws.Range("B24").Value = 'Avalanche'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.32'
$ws.Range("E24").Value = '  +2.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.75'
$ws.Range("E25").Value = '  -0.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.15'
$ws.Range("E26").Value = '  +1.94%  '
$ws.Range("E27").Value = '  +1.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.04'
$ws.Range("E28").Value = '  +3.55%  '
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0500'
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("E31").Value = '  +2.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.32'
$ws.Range("E32").Value = '  +2.62%  '
$ws.Range("D33").Value = '1.481.89'
$ws.Range("E33").Value = '  -2.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.12'
$ws.Range("E34").Value = '  +5.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.62'
$ws.Range("E35").Value = '  +5.90%  '
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.580'
$ws.Range("E37").Value = '  +2.61%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.904'
$ws.Range("E38").Value = '  +9.09%  '
$ws.Range("E39").Value = '  +2.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.07'
$ws.Range("E40").Value = '  +2.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '67.02'
$ws.Range("E42").Value = '  +9.16%  '
$ws.Range("E43").Value = '  +3.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.988'
$ws.Range("E44").Value = '  +8.09%  '
$ws.Range("D45").Value = '1.823.31'
$ws.Range("E45").Value = '  +3.84%  '
$ws.Range("E46").Value = '  +2.16%  '
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("E48").Value = '  +2.59%  '
$ws.Range("E49").Value = '  +5.53%  '
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.72'
$ws.Range("E51").Value = '  +3.01%  '
